$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new "Shipping" column between Quantity (C) and Total (old D)
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "Shipping"

# Arduino Due row: unit price, shipping cost, total formula, and purchase link
$ws.Range("B2").Value = 47.5
$ws.Range("D2").Value = 8.99
$ws.Range("E2").Formula = "=(B2*C2)+D2"
$ws.Range("F2").Value = "http://www.amazon.com/Arduino-A000062-Due/dp/B00A6C3JN2/ref=sr_1_1?ie=UTF8&qid=1422406211&sr=8-1&keywords=arduino+due&pebp=1422406212053&peasin=B00A6C3JN2"

# MCP2551 transceivers row: unit price, shipping cost, total formula, and purchase link
$ws.Range("B3").Value = 1.22
$ws.Range("D3").Value = 5.23
$ws.Range("E3").Formula = "=(B3*C3)+D3"
$ws.Range("F3").Value = "http://www.digikey.com/product-search/en?KeyWords=MCP2551-I%2FP-ND%20&WT.z_header=search_go"

# Apply currency formatting to the price columns
$ws.Range("B2:B11").NumberFormat = """$""#,##0.00"
$ws.Range("D2:D9").NumberFormat = """$""#,##0.00"
$ws.Range("E2:E11").NumberFormat = """$""#,##0.00"

# Extra blank rows with running totals
$ws.Range("E4").Formula = "=SUM(B4:D4)"
$ws.Range("E5:E11").Formula = "=SUM(B5:D5)"

$ws.Range("D3").Select()
